# Fix import of files from intermediate folder.
# Files are now identified by exact matches of full name rather than parts
# of file name, so the "_example" suffix that was used for the intermediate
# / test copies must be dropped from the file names actually configured on
# the "Sheet1" lookup sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Drop the "_example" suffix from the five dataset file names ---
$ws1.Range("B2").Value = "AmphibiansReptiles_Capinha-etal2017.xlsx"
$ws1.Range("B3").Value = "GAVIA_taxon_region_list.xlsx"
$ws1.Range("B4").Value = "GRIIS_sTwist_Hanno_Aug92019_resend.xlsx"
$ws1.Range("B5").Value = "GloNAF_taxon_region_list.xlsx"
$ws1.Range("B6").Value = "GlobalAlienSpeciesFirstRecordDatabase_v1.2_withcountries.xlsx"

# --- Column header rename (space -> dot) ---
$ws1.Range("K4").Value = "Resolved.date"

# --- Row height normalised to the sheet default ---
$ws1.Rows.Item(4).RowHeight = 15

# --- View state: move the frozen-pane viewport / active selection ---
$ws1.Activate() | Out-Null
$ws1.Range("K5").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
